$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "expiry policy" value for row 2 (header row already has H1 = "expiry")
$ws.Range("H2").Value = "session"

# Update the existing expiry timestamps in rows 3 and 4
$ws.Range("H3").Value = 1738815249
$ws.Range("H4").Value = 1736828049
